$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet holds one weekly price observation per row for
# "Femacal de La Calera - Ciboulette" (rows 2..291, header in row 1).
# A new weekly observation is inserted as the new row 71, which pushes
# every existing row from 71..291 down by one (new row 292 is created).
# Only the "Fecha" (D) and the volume/price block (J,K,L,M,N,P) vary
# from row to row; columns A,B,C,E,F,G,H,I,O,Q,R are constant for every
# data row, so the newly created row 292 can just copy them from row 291.

$firstDataRow = 71
$lastOldDataRow = 291
$newLastDataRow = 292

# Copy constant columns into the brand new last row (292) from the
# previous last row (291) before anything else shifts.
$constCols = 1,2,3,5,6,7,8,9,15,17,18
foreach ($c in $constCols) {
    $ws.Cells.Item($newLastDataRow, $c).Value = $ws.Cells.Item($lastOldDataRow, $c).Value2
}

# Shift the date + volume/price block down by one row, working from the
# bottom up so the source row for each step hasn't been overwritten yet.
for ($r = $newLastDataRow; $r -gt $firstDataRow; $r--) {
    $src = $r - 1

    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($src, 4).Value2
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($src, 10).Value2
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item($src, 11).Value2
    $ws.Cells.Item($r, 12).Value = $ws.Cells.Item($src, 12).Value2
    $ws.Cells.Item($r, 13).Value = $ws.Cells.Item($src, 13).Value2
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($src, 14).Value2
    $ws.Cells.Item($r, 16).Value = $ws.Cells.Item($src, 16).Value2
}

# New weekly observation placed at row 71.
$ws.Cells.Item($firstDataRow, 4).Value = 44648
$ws.Cells.Item($firstDataRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($firstDataRow, 10).Value = 180
$ws.Cells.Item($firstDataRow, 11).Value = 1500
$ws.Cells.Item($firstDataRow, 12).Value = 1500
$ws.Cells.Item($firstDataRow, 13).Value = 1500
$ws.Cells.Item($firstDataRow, 14).Value = "$/docena de atados"
$ws.Cells.Item($firstDataRow, 16).Value = 500
